$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.118.69"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "1.796.93"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "222.87"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "32.15"
$ws.Range("E8").Value = "  -1.73%  "
$ws.Range("D9").Value = "0.284"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "0.0717"
$ws.Range("E10").Value = "  +4.61%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "2.054.50"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.789.20"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "10.74"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "0.630"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "34.099.33"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "68.06"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "246.19"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").Value = "0.0₃0784"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "10.76"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("E23").Value = "  -2.33%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "158.84"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "16.48"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "7.04"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "3.71"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").Value = "3.51"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "1.414.28"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "0.942"
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("D40").Value = "80.16"
$ws.Range("E40").Value = "  -3.26%  "
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  +3.63%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -1.62%  "
$ws.Range("D46").Value = "1.953.50"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  -2.92%  "
$ws.Range("D48").Value = "106.15"
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("D50").Value = "11.88"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("E51").Value = "  +0.05%  "
